# #125 updated viewer libs for editable grid and slimmer config
#
# The "person" sample grid gained a 4th "zipcode" column, and its rows were
# re-synced/reordered against the (now editable) source data, which also
# picked up a corrected "age" for one record (the row whose name is "r":
# 906 -> 25). Reproduce that end state directly against the open workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "zipcode" header in column D.
$ws.Range("D1").Value = "zipcode"

# Final row-by-row content for A2:D11 (age, gender, name, zipcode).
# $null marks rows that have no "age" value (those cells stay blank, just
# like in the original sheet).
$data = @(
    @(25,    "MALE",   "r",                    631),
    @(255,   "FEMALE", "oYBuz",                690),
    @($null, "FEMALE", "xKMgdHyLw",             304),
    @($null, "MALE",   "MPPsYun",               875),
    @(167,   "MALE",   "XTrq",                  973),
    @($null, "FEMALE", "g",                     351),
    @($null, "FEMALE", "ynx",                   525),
    @($null, "FEMALE", "xPFZ",                  921),
    @(541,   "MALE",   "bRHCLEwdglb",           123),
    @($null, "FEMALE", "UVcDVhxpyCziyBSiRasp",  53)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rec = $data[$i]

    if ($null -eq $rec[0]) {
        $ws.Cells.Item($row, 1).Value = ""
    } else {
        $ws.Cells.Item($row, 1).Value = $rec[0]
    }
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
}
